# Add a new "CVSSv2 Score" column (between "Categorization" and "CWE Score")
# to both worksheets, matching the vulnmgmt.xlsx commit.
#
# For each sheet, column E ("CVSSv2 Score") is new. Rows that already had a
# value under the old "CWE Score" column (F) shift that value (and the
# Reported/Description/Fixed columns after it) one column to the right, into
# G/H/I/J, leaving E (and now F) blank for that row. The one row that gets an
# actual CVSSv2 value (Tier 1 / sheet1 row 2) keeps its CWE-Score-and-onward
# columns where they were (F was already blank there) and simply gains a
# value in E.

function Set-TextValue {
    param($cell, [string]$value)
    # Force the cell to store its content as text, even when the text looks
    # like a number (e.g. "6.1", "90.0"), mirroring the original workbook's
    # shared-string ("t=s") cell type instead of letting Excel infer a
    # numeric type. Resetting the style back to Normal afterwards keeps the
    # cell on the workbook's default (unstyled) format.
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = "Normal"
}

function Add-CvssColumn {
    param($ws)

    # New header cell, row 1.
    Set-TextValue $ws.Range("E1") "CVSSv2 Score"
}

function Shift-RowRight {
    # Moves F..I of the given row to G..J (one column to the right),
    # clearing F afterwards. Used for rows whose CWE Score column already
    # held a value that needs to make room for the new CVSSv2 Score column.
    param($ws, [int]$row)

    $colF = $ws.Cells.Item($row, 6)
    $colG = $ws.Cells.Item($row, 7)
    $colH = $ws.Cells.Item($row, 8)
    $colI = $ws.Cells.Item($row, 9)

    $vF = $colF.Text
    $vG = $colG.Text
    $vH = $colH.Text
    $vI = $colI.Text

    $colF.ClearContents()

    Set-TextValue $ws.Cells.Item($row, 7) $vF
    Set-TextValue $ws.Cells.Item($row, 8) $vG
    Set-TextValue $ws.Cells.Item($row, 9) $vH
    Set-TextValue $ws.Cells.Item($row, 10) $vI
}

$wb = $excel.ActiveWorkbook

# ---- Sheet "Tier 1" ----
$ws1 = $wb.Worksheets.Item("Tier 1")
Add-CvssColumn $ws1
# Row 2 (42353 / Infrastructure / Patch) gets an actual CVSSv2 score; its
# CWE Score column (F) was already empty, so nothing else needs to move.
Set-TextValue $ws1.Range("E2") "6.1"
# Row 3 (42354 / Application / SQLi) had a CWE Score (90.0) already filled
# in; shift it (and Reported/Description/Fixed) right to make room.
Shift-RowRight $ws1 3

# ---- Sheet "Tier 2" ----
$ws2 = $wb.Worksheets.Item("Tier 2")
Add-CvssColumn $ws2
Shift-RowRight $ws2 2
Shift-RowRight $ws2 3
Shift-RowRight $ws2 4
